$d = $word.ActiveDocument

# The paragraph that ends with "...--skip-tests=true" is the 2nd paragraph
# of the document. We need to insert two new paragraphs right after it:
#   1. A bold/colored sub-heading:
#        "Agregar un componente a otra misma carpeta de otro componente sin añadir carpeta"
#   2. A plain paragraph with the text "--flat"

$skipTestsPara = $d.Paragraphs.Item(2)
$rng = $skipTestsPara.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()

# The freshly-inserted (empty) paragraph is now paragraph #3.
$headingPara = $d.Paragraphs.Item(3)
$headingRng = $headingPara.Range
$headingRng.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:rPr><w:b/><w:bCs/><w:color w:val='C45911' w:themeColor='accent2' w:themeShade='BF'/><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:t>Agregar un componente a otra misma carpeta de otro componente sin añadir carpeta</w:t></w:r></w:p>")

# Insert the second new paragraph right after the heading paragraph.
$headingPara2 = $d.Paragraphs.Item(3)
$rng2 = $headingPara2.Range
$rng2.Collapse(0)
$rng2.InsertParagraphAfter()

$flatPara = $d.Paragraphs.Item(4)
$flatRng = $flatPara.Range
$flatRng.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>--flat</w:t></w:r></w:p>")
